$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I10").Value = 0.5
$ws.Range("I11").Value = 0.5
$ws.Range("I23").Value = 0.5
